# Generate Report for Handoff
#
# Two files that were "Handed back" are now "Ready for handoff" again:
#   d0db8ceb-f17b-4d68-9fe8-e79e3b31081c.md  -> 73c39128-3f88-4cbe-9850-3efc160b9f93.md
#   daa328e7-dc06-43cc-9aa3-e2af1d8ee5db.md  -> ffff8ab0a22f-3722-4d45-85b8-5fc55049da2b.md
#
# The Overview sheet gets the new file names / status / handoff date, and the
# per-locale sheets (zh-cn, de-de) lose their "Latest Target File" /
# "Latest Handback File" columns (F/G) - those are only populated once a
# handback has happened - and pick up the new handoff file name, handoff
# datetime and a reset ("zero") handback datetime.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "73c39128-3f88-4cbe-9850-3efc160b9f93.md"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-57-20 16:57:27"

$ws1.Range("A3").Value = "ffff8ab0a22f-3722-4d45-85b8-5fc55049da2b.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-57-20 16:57:27"

# Rebuild the hyperlinks (display text only - the link targets are left
# pointing at whatever relationship id ends up in that slot, matching the
# source report).
$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0ea90d921043aa7b9938754b5efe4064c08a0d21/e2e/d0db8ceb-f17b-4d68-9fe8-e79e3b31081c.md", "", "", "73c39128-3f88-4cbe-9850-3efc160b9f93.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0ea90d921043aa7b9938754b5efe4064c08a0d21/e2e/daa328e7-dc06-43cc-9aa3-e2af1d8ee5db.md", "", "", "ffff8ab0a22f-3722-4d45-85b8-5fc55049da2b.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "73c39128-3f88-4cbe-9850-3efc160b9f93.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "73c39128-3f88-4cbe-9850-3efc160b9f93.9938c5df5a8f32b29a7ede5bc650f7b859d603c6.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-20 16:57:24"
$ws2.Range("F2:G2").Clear()
$ws2.Range("H2").Value = "0001-01-01 00:00:00"
$ws2.Range("I2").Value = "Include"

$ws2.Range("A3").Value = "ffff8ab0a22f-3722-4d45-85b8-5fc55049da2b.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "73c39128-3f88-4cbe-9850-3efc160b9f93.9938c5df5a8f32b29a7ede5bc650f7b859d603c6.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-20 16:57:24"
$ws2.Range("F3:G3").Clear()
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("I3").Value = "Include"

$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0ea90d921043aa7b9938754b5efe4064c08a0d21/e2e/d0db8ceb-f17b-4d68-9fe8-e79e3b31081c.md", "", "", "73c39128-3f88-4cbe-9850-3efc160b9f93.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/0ea90d921043aa7b9938754b5efe4064c08a0d21/e2e/d0db8ceb-f17b-4d68-9fe8-e79e3b31081c.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7cb5e46404c9df5f6fe5c3facc2cb05165047064/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d0db8ceb-f17b-4d68-9fe8-e79e3b31081c.a60b1b2bc3e9fb83900eedce422449d082d731df.zh-cn.xlf", "", "", "73c39128-3f88-4cbe-9850-3efc160b9f93.9938c5df5a8f32b29a7ede5bc650f7b859d603c6.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fa432443d6a82992076716066a5af69989b30189/e2e/d0db8ceb-f17b-4d68-9fe8-e79e3b31081c.md", "", "", "ffff8ab0a22f-3722-4d45-85b8-5fc55049da2b.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9e131d43b423fac57c832165944126d06cb5a22c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d0db8ceb-f17b-4d68-9fe8-e79e3b31081c.a60b1b2bc3e9fb83900eedce422449d082d731df.zh-cn.xlf", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTest/oltest/blob/0ea90d921043aa7b9938754b5efe4064c08a0d21/e2e/daa328e7-dc06-43cc-9aa3-e2af1d8ee5db.md", "", "", "73c39128-3f88-4cbe-9850-3efc160b9f93.9938c5df5a8f32b29a7ede5bc650f7b859d603c6.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "73c39128-3f88-4cbe-9850-3efc160b9f93.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "73c39128-3f88-4cbe-9850-3efc160b9f93.9938c5df5a8f32b29a7ede5bc650f7b859d603c6.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-20 16:57:27"
$ws3.Range("F2:G2").Clear()
$ws3.Range("H2").Value = "0001-01-01 00:00:00"
$ws3.Range("I2").Value = "Include"

$ws3.Range("A3").Value = "ffff8ab0a22f-3722-4d45-85b8-5fc55049da2b.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "73c39128-3f88-4cbe-9850-3efc160b9f93.9938c5df5a8f32b29a7ede5bc650f7b859d603c6.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-20 16:57:27"
$ws3.Range("F3:G3").Clear()
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("I3").Value = "Include"

$ws3.Cells.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0ea90d921043aa7b9938754b5efe4064c08a0d21/e2e/d0db8ceb-f17b-4d68-9fe8-e79e3b31081c.md", "", "", "73c39128-3f88-4cbe-9850-3efc160b9f93.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/0ea90d921043aa7b9938754b5efe4064c08a0d21/e2e/d0db8ceb-f17b-4d68-9fe8-e79e3b31081c.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/77d6f964407be5c54031bd907dd0b3e069ac6fb4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d0db8ceb-f17b-4d68-9fe8-e79e3b31081c.a60b1b2bc3e9fb83900eedce422449d082d731df.de-de.xlf", "", "", "73c39128-3f88-4cbe-9850-3efc160b9f93.9938c5df5a8f32b29a7ede5bc650f7b859d603c6.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b31330ee53457bc0a0282dc067d0ca1d1dc0bd8a/e2e/d0db8ceb-f17b-4d68-9fe8-e79e3b31081c.md", "", "", "ffff8ab0a22f-3722-4d45-85b8-5fc55049da2b.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/249d072d1b8c1dce69191fcff5d84d1374144b53/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d0db8ceb-f17b-4d68-9fe8-e79e3b31081c.a60b1b2bc3e9fb83900eedce422449d082d731df.de-de.xlf", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTest/oltest/blob/0ea90d921043aa7b9938754b5efe4064c08a0d21/e2e/daa328e7-dc06-43cc-9aa3-e2af1d8ee5db.md", "", "", "73c39128-3f88-4cbe-9850-3efc160b9f93.9938c5df5a8f32b29a7ede5bc650f7b859d603c6.de-de.xlf") | Out-Null

$wb.Save()
